# Apply weekly fruit/vegetable update:
# Insert two new daily records (rows) into the Tomate (Chillán) dataset,
# right before the existing row that used to be row 379 (date 2021-05-07 / serial 44312),
# pushing the remaining rows down by two positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at the top of the block (rows 379 and 380),
# shifting the existing rows 379-400 down to become rows 381-402.
$ws.Rows("379:380").Insert()

# --- New row 379 ---
$ws.Range("A379").Value2 = 7
$ws.Range("B379").Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C379").Value2 = "Ñuble"
$ws.Range("D379").Value2 = 44610
$ws.Range("E379").Value2 = 16
$ws.Range("F379").Value2 = 100112020
$ws.Range("G379").Value2 = "Tomate"
$ws.Range("H379").Value2 = "Larga vida"
$ws.Range("I379").Value2 = "Primera"
$ws.Range("J379").Value2 = 800
$ws.Range("K379").Value2 = 6000
$ws.Range("L379").Value2 = 6500
$ws.Range("M379").Value2 = 6250
$ws.Range("N379").Value2 = "`$/caja 15 kilos"
$ws.Range("O379").Value2 = "Región del Maule"
$ws.Range("P379").Value2 = 417
$ws.Range("Q379").Value2 = 15
$ws.Range("R379").Value2 = "Hortaliza"

# --- New row 380 ---
$ws.Range("A380").Value2 = 7
$ws.Range("B380").Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C380").Value2 = "Ñuble"
$ws.Range("D380").Value2 = 44610
$ws.Range("E380").Value2 = 16
$ws.Range("F380").Value2 = 100112020
$ws.Range("G380").Value2 = "Tomate"
$ws.Range("H380").Value2 = "Larga vida"
$ws.Range("I380").Value2 = "Segunda"
$ws.Range("J380").Value2 = 400
$ws.Range("K380").Value2 = 5000
$ws.Range("L380").Value2 = 5500
$ws.Range("M380").Value2 = 5250
$ws.Range("N380").Value2 = "`$/caja 15 kilos"
$ws.Range("O380").Value2 = "Región del Maule"
$ws.Range("P380").Value2 = 350
$ws.Range("Q380").Value2 = 15
$ws.Range("R380").Value2 = "Hortaliza"

# Keep the date-formatted style on D379/D380 consistent with the rest of column D.
$ws.Range("D379:D380").NumberFormat = $ws.Range("D381").NumberFormat
